# Updated capital structure database
# Applies the diff: row 2 (Malta / "1") gets refreshed figures and becomes "2",
# a brand-new row for "Santumas Shareholdings plc (MTSE:STS)" is inserted as
# the new row 3, and the former row 3 ("Brait SE (JSE:BAT)") shifts down to
# row 4 with refreshed figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Update existing row 2 (Malta / company id) in place with refreshed data
# ---------------------------------------------------------------------------
$ws.Cells.Item(2,1).Value = "Malta"

# B2 ("1" -> "2") must remain stored as text, not get auto-converted to a number
$ws.Cells.Item(2,2).NumberFormat = "@"
$ws.Cells.Item(2,2).Value = "2"
$ws.Cells.Item(2,2).Style = "Normal"

$ws.Cells.Item(2,3).Value = "Investments & Asset Management"

$ws.Cells.Item(2,7).Value = 0.0007089204182826845
$ws.Cells.Item(2,8).Value = 0.0007089204182826845
$ws.Cells.Item(2,9).Value = 1.020865040011783
$ws.Cells.Item(2,10).Value = 1.020865040011783
$ws.Cells.Item(2,11).Value = -1065.66
$ws.Cells.Item(2,12).Value = 1.046354754774412
$ws.Cells.Item(2,13).Value = 0
$ws.Cells.Item(2,14).Value = 0
$ws.Cells.Item(2,15).Value = -0
$ws.Cells.Item(2,16).Value = 0
$ws.Cells.Item(2,17).Value = 0
$ws.Cells.Item(2,18).Value = -0
$ws.Cells.Item(2,19).Value = 0
$ws.Cells.Item(2,20).ClearContents()
$ws.Cells.Item(2,21).Value = 12.333
$ws.Cells.Item(2,22).Value = 0.03539896670493685
$ws.Cells.Item(2,23).Value = -0.5153716287001916
$ws.Cells.Item(2,24).Value = 0.04962139047014914
$ws.Cells.Item(2,25).Value = -0.5649930191703407
$ws.Cells.Item(2,26).Value = -0.5124000871399628
$ws.Cells.Item(2,27).Value = -0.3291813778731087
$ws.Cells.Item(2,28).Value = 0.04017937176197776
$ws.Cells.Item(2,29).Value = -0.3693607496350865
$ws.Cells.Item(2,30).Value = 334.7
$ws.Cells.Item(2,31).Value = 0
$ws.Cells.Item(2,32).Value = 334.7
$ws.Cells.Item(2,33).Value = 322.367
$ws.Cells.Item(2,34).Value = 0.4899721856243596
$ws.Cells.Item(2,35).Value = 0.3500313741895001
$ws.Cells.Item(2,36).Value = 0.4805946028948949
$ws.Cells.Item(2,37).Value = 0.3415385854151061
$ws.Cells.Item(2,38).Value = 76.3
$ws.Cells.Item(2,39).Value = 76.3
$ws.Cells.Item(2,40).Value = -198.0473372781065
$ws.Cells.Item(2,41).Value = -13.62647444298821
$ws.Cells.Item(2,42).Value = -190.7497041420118
$ws.Cells.Item(2,43).Value = -13.62647444298821

# ---------------------------------------------------------------------------
# 2) Insert a brand-new row 3 (pushes the former row 3 "Brait SE" down to 4)
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).Insert()

# ---------------------------------------------------------------------------
# 3) Populate the new row 3 - Santumas Shareholdings plc (MTSE:STS)
# ---------------------------------------------------------------------------
$ws.Cells.Item(3,1).Value = "Malta"
$ws.Cells.Item(3,2).Value = "Santumas Shareholdings plc (MTSE:STS)"
$ws.Cells.Item(3,3).Value = "Investments & Asset Management"

$ws.Cells.Item(3,7).Value = 0.4658064516129032
$ws.Cells.Item(3,8).Value = 0.4658064516129032
$ws.Cells.Item(3,9).Value = 1.096774193548387
$ws.Cells.Item(3,10).Value = 1.096774193548387
$ws.Cells.Item(3,11).Value = -1.76
$ws.Cells.Item(3,12).Value = 1.135483870967742
$ws.Cells.Item(3,13).Value = -0
$ws.Cells.Item(3,14).Value = -0
$ws.Cells.Item(3,15).Value = 0
$ws.Cells.Item(3,16).Value = -0
$ws.Cells.Item(3,17).Value = -0
$ws.Cells.Item(3,18).Value = 0
$ws.Cells.Item(3,19).Value = 0
$ws.Cells.Item(3,21).Value = 0.633
$ws.Cells.Item(3,22).Value = 0.05104838709677419
$ws.Cells.Item(3,23).Value = -0.1323308270676692
$ws.Cells.Item(3,24).Value = 0.03976014512418032
$ws.Cells.Item(3,25).Value = -0.1720909721918495
$ws.Cells.Item(3,26).Value = -0.1210275630514562
$ws.Cells.Item(3,27).Value = -0.1327399078628875
$ws.Cells.Item(3,28).Value = 0.03976014512418032
$ws.Cells.Item(3,29).Value = -0.1725000529870678
$ws.Cells.Item(3,30).Value = 0
$ws.Cells.Item(3,31).Value = 0
$ws.Cells.Item(3,32).Value = 0
$ws.Cells.Item(3,33).Value = -0.633
$ws.Cells.Item(3,34).Value = 0
$ws.Cells.Item(3,35).Value = 0
$ws.Cells.Item(3,36).Value = -0.05379451007053625
$ws.Cells.Item(3,37).Value = -0.0552018836661725
$ws.Cells.Item(3,38).Value = 0
$ws.Cells.Item(3,39).Value = 0
$ws.Cells.Item(3,40).Value = -0
$ws.Cells.Item(3,42).Value = 0.3745562130177515

# ---------------------------------------------------------------------------
# 4) Update row 4 (former row 3, "Brait SE (JSE:BAT)") with refreshed figures
# ---------------------------------------------------------------------------
$ws.Cells.Item(4,9).Value = 1.020749336217917
$ws.Cells.Item(4,10).Value = 1.020749336217917
$ws.Cells.Item(4,11).Value = -1063.9
$ws.Cells.Item(4,12).Value = 1.046218900580195
$ws.Cells.Item(4,13).Value = -0
$ws.Cells.Item(4,14).Value = -0
$ws.Cells.Item(4,15).Value = 0
$ws.Cells.Item(4,19).Value = 0
$ws.Cells.Item(4,20).ClearContents()
$ws.Cells.Item(4,21).Value = 11.7
$ws.Cells.Item(4,22).Value = 0.03482142857142857
$ws.Cells.Item(4,23).Value = -0.8984124303327141
$ws.Cells.Item(4,24).Value = 0.05948263581611797
$ws.Cells.Item(4,25).Value = -0.9578950661488321
$ws.Cells.Item(4,26).Value = -0.5149382215920599
$ws.Cells.Item(4,27).Value = -0.52562284788333
$ws.Cells.Item(4,28).Value = 0.04059859839977521
$ws.Cells.Item(4,29).Value = -0.5662214462831052
$ws.Cells.Item(4,30).Value = 334.7
$ws.Cells.Item(4,32).Value = 334.7
$ws.Cells.Item(4,33).Value = 323
$ws.Cells.Item(4,34).Value = 0.4990308632771731
$ws.Cells.Item(4,35).Value = 0.3545175299226777
$ws.Cells.Item(4,36).Value = 0.4901365705614568
$ws.Cells.Item(4,37).Value = 0.3464178464178465
$ws.Cells.Item(4,38).Value = 76.3
$ws.Cells.Item(4,39).Value = 76.3
$ws.Cells.Item(4,41).Value = -13.60419397116645
$ws.Cells.Item(4,43).Value = -13.60419397116645
